$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New labelled settings block (rows 14, 16, 18) ---
# Introduce new shared strings in the exact order they first appear so the
# rebuilt sharedStrings table lines up with the target workbook.
$ws.Range("H14").Value = "PenWidth"
$ws.Range("I14").Value = "text"

$ws.Range("H16").Value = "Tunels count"
$ws.Range("H18").Value = "Spheres count"

$ws.Range("K14").Value = "reverse map"
$ws.Range("K16").Value = "tunels"
$ws.Range("K18").Value = "spheres"

# Re-used strings (already present elsewhere in the sheet / sst).
$ws.Range("I16").Value = "text"
$ws.Range("I18").Value = "text"
$ws.Range("O16").Value = "Save Conf"
$ws.Range("O18").Value = "Load Conf"

# --- Remove the old "Map/Game/Weap/Team Set" labels from row 8 & 10 ---
# and give those whole rows (F:P) the new yellow-highlight style.
$rowsToHighlight = 5,6,7,8,9,10,11,12
foreach ($r in $rowsToHighlight) {
    $rng = $ws.Range("F" + $r + ":P" + $r)
    $rng.ClearContents()
    $rng.Interior.Color = 65535
}

# --- Drop the old numeric demo values from row 16 that aren't part of the
#     new layout, and the "Start" label in row 18 ---
$ws.Range("J16").ClearContents()
$ws.Range("L16").ClearContents()
$ws.Range("N16").ClearContents()
$ws.Range("P18").ClearContents()

# --- Update the view: scroll so column D is left-most and select K18 ---
$ws.Range("K18").Select()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
